$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.994.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.908.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.87%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4806'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +0.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07356'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9321'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.77'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07760'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.926.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.637'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008828'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.022.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.161.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.917'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.124'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08932'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.309'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.259'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7743'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.677'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.622'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02051'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05301'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.001'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5482'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.017'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1524'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.463'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4830'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '108.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.650'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("E51").Value = '  -0.18%  '
